$wb = $excel.ActiveWorkbook

# --- Fix LaTeX strings in params_deep (strip stray leading backslash-escaping) ---
$wsDeep = $wb.Worksheets.Item("params_deep")
$wsDeep.Range("B2").Value  = "\theta_{k}"
$wsDeep.Range("B3").Value  = "\theta_{s}"
$wsDeep.Range("B4").Value  = "\gamma"
$wsDeep.Range("B5").Value  = "\delta_{k}"
$wsDeep.Range("B6").Value  = "\delta_{s}"
$wsDeep.Range("B7").Value  = "\bar{r}"
$wsDeep.Range("B12").Value = "\bar{d}"
$wsDeep.Range("B13").Value = "\varkappa"
$wsDeep.Range("B14").Value = "\rho_{1}"
$wsDeep.Range("B15").Value = "\rho_{2}"
$wsDeep.Range("B16").Value = "\psi"
$wsDeep.Range("B18").Value = "\omega"
$wsDeep.Range("B21").Value = "\text{itermax}"

# --- Fix latex string in params_init ---
$wsInit = $wb.Worksheets.Item("params_init")
$wsInit.Range("B11").Value = "govexp_{0}"

# --- Switch the active sheet / selection: equations -> params_init, cell B12 ---
$wsInit.Activate()
$wsInit.Range("B12").Select()
